# Rename the "alpha" sheet to "summary" and make it the active/selected
# sheet (previously "E1_familiar" was the selected tab), moving the
# selection on that sheet from O27 to L22.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("alpha")
$ws.Name = "summary"

$ws.Activate()
$ws.Range("L22").Select()
